$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (row 2 and row 3) need to have their values swapped
# for columns D, J, K, L, M, N, O, P, Q (the other columns are identical
# between the two rows so no visible change occurs there).

$columns = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

foreach ($col in $columns) {
    $cellRow2 = $ws.Range("$col`2")
    $cellRow3 = $ws.Range("$col`3")

    $valRow2 = $cellRow2.Value2
    $valRow3 = $cellRow3.Value2

    $cellRow2.Value = $valRow3
    $cellRow3.Value = $valRow2
}
